$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Copy formatting (number format, font, borders, alignment) from column F
# into the new column G, so the added column matches the existing styling.
$ws.Range("F1:F30").Copy() | Out-Null
$ws.Range("G1:G30").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Give column G its own width (stored width 17).
$ws.Columns.Item(7).ColumnWidth = 16.14

# Header
$ws.Range("G1").Value = "PRESUPUESTO"

# Data rows 2-29 -> 0
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 7).Value = 0
}

# Total row 30 -> 0
$ws.Range("G30").Value = 0
